# Insert a new data row at row 15 (pushes existing rows 15..113 down to 16..114)
# and populate it with the new week's record, matching the author's commit:
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value2 = 44462
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 100114001
$ws.Range("G15").Value = "Papa"
$ws.Range("H15").Value = "Asterix"
$ws.Range("I15").Value = "1a (guarda)"
$ws.Range("J15").Value = 2000
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9500
$ws.Range("M15").Value = 9250
$ws.Range("N15").Value = "`$/saco 25 kilos"
$ws.Range("O15").Value = "Provincia de Arauco"
$ws.Range("P15").Value = 370
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
